# Add manufacturer / authorised-rep account details for priya (row 21 & 22 of
# the "TestEnv" sheet) and make "TestEnv" the active / selected sheet.

$wb = $excel.ActiveWorkbook

$wsTestEnv = $wb.Worksheets.Item("TestEnv")

# New username values for Priya's manufacturer / authorised-rep rows.
$wsTestEnv.Range("A21").Value = "Manufacturer1681H8_PG"
$wsTestEnv.Range("A22").Value = "AuthorisedRep1681H9_PG"

# Move the active tab / selection from PreProdEnv to TestEnv, with the
# cursor resting on A22 (the row just updated).
$wsTestEnv.Select()
$wsTestEnv.Range("A22").Select()

$wb.Save()
